$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.6591047123922635
$ws.Range("C2").Value2 = 0.06979013589212002
$ws.Range("D2").Value2 = 0.07649305032681841
$ws.Range("E2").Value2 = 0.1128240276633434
$ws.Range("F2").Value2 = 1.807478639038337
$ws.Range("K2").Value2 = 0.6515584640497707
$ws.Range("M2").Value2 = 0.3084085407963784
$ws.Range("N2").Value2 = 2.635021442869899
$ws.Range("B3").Value2 = 0.614129082425336
$ws.Range("C3").Value2 = 0.06204465075281007
$ws.Range("D3").Value2 = 0.07671788445026362
$ws.Range("E3").Value2 = 0.1034124664790923
$ws.Range("F3").Value2 = 1.775948807909671
$ws.Range("K3").Value2 = 0.601413058193998
$ws.Range("M3").Value2 = 0.2838194162092762
$ws.Range("N3").Value2 = 2.632067726335137
$ws.Range("B4").Value2 = 0.5869108781930663
$ws.Range("C4").Value2 = 0.05731793280592967
$ws.Range("D4").Value2 = 0.07685709186232259
$ws.Range("E4").Value2 = 0.09770089020241102
$ws.Range("F4").Value2 = 1.757449805792263
$ws.Range("K4").Value2 = 0.5710057226789047
$ws.Range("M4").Value2 = 0.268903493606075
$ws.Range("N4").Value2 = 2.630894410618978
$ws.Range("B5").Value2 = 0.5759186639841687
$ws.Range("C5").Value2 = 0.0553989287663228
$ws.Range("D5").Value2 = 0.07691410927462616
$ws.Range("E5").Value2 = 0.09538999018739958
$ws.Range("F5").Value2 = 1.750127065168215
$ws.Range("K5").Value2 = 0.5587100363281081
$ws.Range("M5").Value2 = 0.2628704109427957
$ws.Range("N5").Value2 = 2.630576868961796
$ws.Range("B6").Value2 = 0.5740994085054467
$ws.Range("C6").Value2 = 0.05508070898169137
$ws.Range("D6").Value2 = 0.0769235943878801
$ws.Range("E6").Value2 = 0.09500726244970537
$ws.Range("F6").Value2 = 1.748924142906446
$ws.Range("K6").Value2 = 0.5566741022614679
$ws.Range("M6").Value2 = 0.2618713433420083
$ws.Range("N6").Value2 = 2.63053382900749
$ws.Range("B7").Value2 = 0.586762231313827
$ws.Range("C7").Value2 = 0.05729202357890983
$ws.Range("D7").Value2 = 0.07685785965099257
$ws.Range("E7").Value2 = 0.09766965768800873
$ws.Range("F7").Value2 = 1.757350175994489
$ws.Range("K7").Value2 = 0.5708395125881509
$ws.Range("M7").Value2 = 0.2688219466259767
$ws.Range("N7").Value2 = 2.630889478464795
$ws.Range("B8").Value2 = 0.6435144523369729
$ws.Range("C8").Value2 = 0.06711338372225839
$ws.Range("D8").Value2 = 0.07657033150621739
$ws.Range("E8").Value2 = 0.1095648137322627
$ws.Range("F8").Value2 = 1.79642818471433
$ws.Range("K8").Value2 = 0.6341886674001671
$ws.Range("M8").Value2 = 0.2998921696833179
$ws.Range("N8").Value2 = 2.633869860455789
$ws.Range("B9").Value2 = 0.7579806014014423
$ws.Range("C9").Value2 = 0.08661038358320639
$ws.Range("D9").Value2 = 0.07601575579679576
$ws.Range("E9").Value2 = 0.1334376943438045
$ws.Range("F9").Value2 = 1.879920921771287
$ws.Range("K9").Value2 = 0.7614801266465463
$ws.Range("M9").Value2 = 0.362288387711196
$ws.Range("N9").Value2 = 2.644815051827678
$ws.Range("B10").Value2 = 0.8440589803890362
$ws.Range("C10").Value2 = 0.1010905547005052
$ws.Range("D10").Value2 = 0.07561403926779153
$ws.Range("E10").Value2 = 0.1513309495337651
$ws.Range("F10").Value2 = 1.945500955783643
$ws.Range("K10").Value2 = 0.8569235279091743
$ws.Range("M10").Value2 = 0.4090641847685745
$ws.Range("N10").Value2 = 2.655996909133009
$ws.Range("B11").Value2 = 0.8836584875021458
$ws.Range("C11").Value2 = 0.1077142584298372
$ws.Range("D11").Value2 = 0.0754325522274506
$ws.Range("E11").Value2 = 0.1595524870687512
$ws.Range("F11").Value2 = 1.9762678306559
$ws.Range("K11").Value2 = 0.9007731225271414
$ws.Range("M11").Value2 = 0.4305550452683207
$ws.Range("N11").Value2 = 2.661772659867182
$ws.Range("B12").Value2 = 0.8987179914253147
$ws.Range("C12").Value2 = 0.1102279179746688
$ws.Range("D12").Value2 = 0.07536401123986636
$ws.Range("E12").Value2 = 0.1626778713439379
$ws.Range("F12").Value2 = 1.988053640294822
$ws.Range("K12").Value2 = 0.9174407350374736
$ws.Range("M12").Value2 = 0.4387242121703565
$ws.Range("N12").Value2 = 2.664059396719466
$ws.Range("B13").Value2 = 0.8954718025056536
$ws.Range("C13").Value2 = 0.1096863136490356
$ws.Range("D13").Value2 = 0.0753787645142836
$ws.Range("E13").Value2 = 0.1620042224059333
$ws.Range("F13").Value2 = 1.985509336713733
$ws.Range("K13").Value2 = 0.9138482669618782
$ws.Range("M13").Value2 = 0.4369634469569235
$ws.Range("N13").Value2 = 2.663562469928792
$ws.Range("B14").Value2 = 0.8848961559833697
$ws.Range("C14").Value2 = 0.1079209495268572
$ws.Range("D14").Value2 = 0.07542690960509457
$ws.Range("E14").Value2 = 0.1598093707382446
$ws.Range("F14").Value2 = 1.977234745268106
$ws.Range("K14").Value2 = 0.90214311642805
$ws.Range("M14").Value2 = 0.4312265024039448
$ws.Range("N14").Value2 = 2.661958792369745
$ws.Range("B15").Value2 = 0.8784266197642978
$ws.Range("C15").Value2 = 0.1068403214294733
$ws.Range("D15").Value2 = 0.07545642396535879
$ws.Range("E15").Value2 = 0.1584665411665611
$ws.Range("F15").Value2 = 1.972183929869118
$ws.Range("K15").Value2 = 0.8949815604990476
$ws.Range("M15").Value2 = 0.4277165164709515
$ws.Range("N15").Value2 = 2.660989477408748
$ws.Range("B16").Value2 = 0.8414799938524311
$ws.Range("C16").Value2 = 0.1006584288195143
$ws.Range("D16").Value2 = 0.07562592544728108
$ws.Range("E16").Value2 = 0.150795324681738
$ws.Range("F16").Value2 = 1.94350912829384
$ws.Range("K16").Value2 = 0.8540665972281829
$ws.Range("M16").Value2 = 0.4076640225996755
$ws.Range("N16").Value2 = 2.65563335849825
$ws.Range("B17").Value2 = 0.8189279777695617
$ws.Range("C17").Value2 = 0.09687551391684224
$ws.Range("D17").Value2 = 0.07573023327249384
$ws.Range("E17").Value2 = 0.146110448795163
$ws.Range("F17").Value2 = 1.92615784811494
$ws.Range("K17").Value2 = 0.8290776494208103
$ws.Range("M17").Value2 = 0.3954172098195556
$ws.Range("N17").Value2 = 2.652524392396927
$ws.Range("B18").Value2 = 0.8059982311523584
$ws.Range("C18").Value2 = 0.09470311771536899
$ws.Range("D18").Value2 = 0.07579034634717097
$ws.Range("E18").Value2 = 0.1434235137650859
$ws.Range("F18").Value2 = 1.916265736855195
$ws.Range("K18").Value2 = 0.8147452838876177
$ws.Range("M18").Value2 = 0.3883931201024495
$ws.Range("N18").Value2 = 2.650801025497586
$ws.Range("B19").Value2 = 0.8016275678138527
$ws.Range("C19").Value2 = 0.09396816755074155
$ws.Range("D19").Value2 = 0.07581071975624898
$ws.Range("E19").Value2 = 0.1425150738912251
$ws.Range("F19").Value2 = 1.912931510505743
$ws.Range("K19").Value2 = 0.8098995485280227
$ws.Range("M19").Value2 = 0.3860182954256786
$ws.Range("N19").Value2 = 2.650228642281348
$ws.Range("B20").Value2 = 0.8213243720046819
$ws.Range("C20").Value2 = 0.09727785497747732
$ws.Range("D20").Value2 = 0.07571911729255643
$ws.Range("E20").Value2 = 0.1466083646194534
$ws.Range("F20").Value2 = 1.92799582042332
$ws.Range("K20").Value2 = 0.8317335575212326
$ws.Range("M20").Value2 = 0.3967188351723507
$ws.Range("N20").Value2 = 2.652848634129384
$ws.Range("B21").Value2 = 0.8880007381941368
$ws.Range("C21").Value2 = 0.1084393318741661
$ws.Range("D21").Value2 = 0.07541276319528301
$ws.Range("E21").Value2 = 0.1604537219679969
$ws.Range("F21").Value2 = 1.979661521970257
$ws.Range("K21").Value2 = 0.9055794956545924
$ws.Range("M21").Value2 = 0.4329107354226096
$ws.Range("N21").Value2 = 2.662427124620194
$ws.Range("B22").Value2 = 0.9319510374988909
$ws.Range("C22").Value2 = 0.1157655844965575
$ws.Range("D22").Value2 = 0.07521361868617849
$ws.Range("E22").Value2 = 0.1695729386689493
$ws.Range("F22").Value2 = 2.014215633273835
$ws.Range("K22").Value2 = 0.9542079734534923
$ws.Range("M22").Value2 = 0.4567454280979177
$ws.Range("N22").Value2 = 2.66926791746593
$ws.Range("B23").Value2 = 0.9084596088638932
$ws.Range("C23").Value2 = 0.1118524879332767
$ws.Range("D23").Value2 = 0.0753198061809357
$ws.Range("E23").Value2 = 0.1646992916983478
$ws.Range("F23").Value2 = 1.995701146475
$ws.Range("K23").Value2 = 0.9282203552013755
$ws.Range("M23").Value2 = 0.4440076534438617
$ws.Range("N23").Value2 = 2.665563557624296
$ws.Range("B24").Value2 = 0.8202408509839074
$ws.Range("C24").Value2 = 0.09709594895801388
$ws.Range("D24").Value2 = 0.07572414237468994
$ws.Range("E24").Value2 = 0.1463832367712214
$ws.Range("F24").Value2 = 1.927164613501532
$ws.Range("K24").Value2 = 0.8305327152614552
$ws.Range("M24").Value2 = 0.3961303181694902
$ws.Range("N24").Value2 = 2.652701845052064
$ws.Range("B25").Value2 = 0.726669988826643
$ws.Range("C25").Value2 = 0.08130938494248596
$ws.Range("D25").Value2 = 0.07616478366078372
$ws.Range("E25").Value2 = 0.1269188124458438
$ws.Range("F25").Value2 = 1.856593735678018
$ws.Range("K25").Value2 = 0.7267107997775213
$ws.Range("M25").Value2 = 0.3452474711761369
$ws.Range("N25").Value2 = 2.641304851207053
